$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 27 de Agosto de 2020 a las 04:21'
$ws.Cells.Item(4, 2).Value = 6000365
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 3313861
$ws.Cells.Item(4, 5).Value = 2502851
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 183653
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(10, 1).Value = 'Mexico'
$ws.Cells.Item(10, 2).Value = 573888
$ws.Cells.Item(10, 3).Value = 5267
$ws.Cells.Item(10, 4).Value = 396758
$ws.Cells.Item(10, 5).Value = 115054
$ws.Cells.Item(10, 7).Value = 626
$ws.Cells.Item(10, 8).Value = 62076
$ws.Cells.Item(11, 1).Value = 'Colombia'
$ws.Cells.Item(11, 2).Value = 572270
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 407121
$ws.Cells.Item(11, 5).Value = 146965
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 18184
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(29, 2).Value = 112094
$ws.Cells.Item(29, 3).Value = 1095
$ws.Cells.Item(29, 4).Value = 50397
$ws.Cells.Item(29, 5).Value = 56971
$ws.Cells.Item(29, 7).Value = 62
$ws.Cells.Item(29, 8).Value = 4726
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(38, 2).Value = 85004
$ws.Cells.Item(38, 3).Value = 8
$ws.Cells.Item(38, 4).Value = 80046
$ws.Cells.Item(38, 5).Value = 324
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(49, 1).Value = 'Honduras'
$ws.Cells.Item(49, 2).Value = 56649
$ws.Cells.Item(49, 3).Value = 772
$ws.Cells.Item(49, 4).Value = 9169
$ws.Cells.Item(49, 5).Value = 45733
$ws.Cells.Item(49, 7).Value = 44
$ws.Cells.Item(49, 8).Value = 1747
$ws.Cells.Item(50, 1).Value = 'Singapur'
$ws.Cells.Item(50, 2).Value = 56495
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 54971
$ws.Cells.Item(50, 5).Value = 1497
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 27
$ws.Cells.Item(51, 1).Value = 'Portugal'
$ws.Cells.Item(51, 2).Value = 56274
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 41184
$ws.Cells.Item(51, 5).Value = 13283
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 1807
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(67, 3).Value = 0
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(68, 3).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(69, 3).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(70, 3).Value = 0
$ws.Cells.Item(71, 3).Value = 0
$ws.Cells.Item(72, 2).Value = 25322
$ws.Cells.Item(72, 3).Value = 118
$ws.Cells.Item(72, 5).Value = 4650
$ws.Cells.Item(72, 7).Value = 23
$ws.Cells.Item(72, 8).Value = 572
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(74, 3).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(75, 3).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(77, 1).Value = 'Corea del Sur'
$ws.Cells.Item(77, 2).Value = 18706
$ws.Cells.Item(77, 3).Value = 441
$ws.Cells.Item(77, 4).Value = 14461
$ws.Cells.Item(77, 5).Value = 3932
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 313
$ws.Cells.Item(78, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(78, 2).Value = 18609
$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(78, 4).Value = 12336
$ws.Cells.Item(78, 5).Value = 5702
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 571
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(82, 1).Value = 'Paraguay'
$ws.Cells.Item(82, 2).Value = 14872
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = 8134
$ws.Cells.Item(82, 5).Value = 6491
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 247
$ws.Cells.Item(83, 1).Value = 'Madagascar'
$ws.Cells.Item(83, 2).Value = 14554
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 13582
$ws.Cells.Item(83, 5).Value = 791
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 181
$ws.Cells.Item(84, 1).Value = 'Libano'
$ws.Cells.Item(84, 2).Value = 14248
$ws.Cells.Item(84, 4).Value = 3955
$ws.Cells.Item(84, 5).Value = 10154
$ws.Cells.Item(84, 8).Value = 139
$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(89, 3).Value = 0
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(91, 3).Value = 0
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(94, 3).Value = 0
$ws.Cells.Item(95, 3).Value = 0
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(100, 3).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(111, 5).Value = 496
$ws.Cells.Item(111, 8).Value = 79
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(117, 3).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(139, 3).Value = 0
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(142, 3).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(147, 2).Value = 1702
$ws.Cells.Item(147, 3).Value = 7
$ws.Cells.Item(147, 4).Value = 1554
$ws.Cells.Item(147, 5).Value = 126
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(166, 3).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(169, 3).Value = 0
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(173, 2).Value = 464
$ws.Cells.Item(173, 3).Value = 33
$ws.Cells.Item(173, 5).Value = 359
$ws.Cells.Item(173, 7).Value = 1
$ws.Cells.Item(173, 8).Value = 3
$ws.Cells.Item(174, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(174, 2).Value = 442
$ws.Cells.Item(174, 3).Value = 24
$ws.Cells.Item(174, 4).Value = 179
$ws.Cells.Item(174, 5).Value = 246
$ws.Cells.Item(174, 8).Value = 17
$ws.Cells.Item(175, 1).Value = 'Burundi'
$ws.Cells.Item(175, 2).Value = 430
$ws.Cells.Item(175, 4).Value = 345
$ws.Cells.Item(175, 5).Value = 84
$ws.Cells.Item(175, 8).Value = 1
$ws.Cells.Item(176, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(176, 2).Value = 419
$ws.Cells.Item(176, 3).Value = 0
$ws.Cells.Item(176, 4).Value = 232
$ws.Cells.Item(176, 5).Value = 183
$ws.Cells.Item(176, 8).Value = 4
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(197, 3).Value = 0
